$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3703749235513688
$ws.Range("D2").Value = 0.3703749235513688

$ws.Range("B3").Value = 0.02109294287051686
$ws.Range("C3").Value = 0.02069750800941871
$ws.Range("D3").Value = 0.02090606863629583

$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01964285593392345
$ws.Range("C4").Value = 0.01943550176427579
$ws.Range("D4").Value = 0.01949951016635842

$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.02987875011818537
$ws.Range("C5").Value = 0.02822797799994848
$ws.Range("D5").Value = 0.02344033290274598
